# Merging prediction and classification into a single repo
# Update classification prediction/error (and one cross-entropy-loss) values
# for the ANN_128nodes ... dropout0.5 results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Control 46
$ws.Range("D2").Value = 2.135435661791767 * [Math]::Pow(10, -5)
$ws.Range("E2").Value = 2.135435661791767 * [Math]::Pow(10, -5)

# Row 3 - Control 28
$ws.Range("D3").Value = 0.9999999903232489
$ws.Range("E3").Value = 0.9999999903232489

# Row 4 - Control 13
$ws.Range("D4").Value = 0.0319742801520539
$ws.Range("E4").Value = 0.0319742801520539

# Row 5 - Control 50
$ws.Range("D5").Value = 5.113401797805511 * [Math]::Pow(10, -45)
$ws.Range("E5").Value = 5.113401797805511 * [Math]::Pow(10, -45)

# Row 6 - Control 51
$ws.Range("D6").Value = 0.8587544168472843
$ws.Range("E6").Value = 0.8587544168472843

# Row 7 - MDD 4
$ws.Range("D7").Value = 0.009001764730336983
$ws.Range("E7").Value = 0.990998235269663

# Row 8 - MDD 32
$ws.Range("D8").Value = 0.9397863203476604
$ws.Range("E8").Value = 0.06021367965233959

# Row 10 - MDD 44
$ws.Range("D10").Value = 0.9991223946686987
$ws.Range("E10").Value = 0.0008776053313013232

# Row 11 - MDD 31
$ws.Range("D11").Value = 0.9999999483468747
$ws.Range("E11").Value = 5.165312533250699 * [Math]::Pow(10, -8)
$ws.Range("F11").Value = 2.521662712097168
